# Masterdata utility update for device_type.xlsx
# Rebuild the sheet to match the new "device_type" masterdata upload format:
# columns: (A) id, (B) lang_code, (C) code, (D) name, (E) descr, (F) is_active

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old 2-row sample data ------------------------------------
$ws.Range("A1:E2").Clear()

# --- New header row (row 1) ---------------------------------------------
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "descr"
$ws.Range("F1").Value = "is_active"

# Header formatting: bold font, thin border all around, centered/top aligned
$headerRange = $ws.Range("B1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows -------------------------------------------------------------
$data = @(
    @(0, "eng", "FRS", "Finger Print Scanner",          "For scanning fingerprints"),
    @(1, "fra", "FRS", "Scanner dempreintes digitales",  "Scannez les empreintes digitales"),
    @(2, "eng", "IRS", "Iris Scanner",                   "For scanning Iris"),
    @(3, "fra", "IRS", "Scanner dIris",                  "Pour scanner liris"),
    @(4, "eng", "CMR", "Camera",                         "For capturing photo"),
    @(5, "fra", "CMR", "Caméra",                          "Pour capturer une photo"),
    @(6, "eng", "SCN", "Document Scanner",               "For scanning documents"),
    @(7, "fra", "SCN", "Scanner de documents",           "Pour numériser des documents"),
    @(8, "eng", "PRT", "Printer",                        "For printing Documents"),
    @(9, "fra", "PRT", "Imprimante",                     "Pour imprimer des documents")
)

$row = 2
foreach ($entry in $data) {
    $idCell = $ws.Cells.Item($row, 1)
    $idCell.Value = $entry[0]
    $idCell.Font.Bold = $true
    $idCell.HorizontalAlignment = -4108
    $idCell.VerticalAlignment = -4160
    $idCell.Borders.LineStyle = 1

    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $true
    $row = $row + 1
}
